# Apply the "mislabelled primary infection site" correction:
#  - refresh several summary-statistic rows with updated source numbers
#  - remove the duplicate "Primary/specific infection site" row (old row 47),
#    which shifts the infection-site rows below it up by one
#  - the row that becomes "Urinary tract" (old row 50 -> new row 49) also
#    picks up corrected figures
#  - the final row (old row 55 -> new row 54) gains an A-column label

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Beta-lactams ---
$ws.Range("D2").Value = 23
$ws.Range("E2").Value = 4022
$ws.Range("F2").Value = 8671
$ws.Range("G2").Value = 1566
$ws.Range("H2").Value = 1404
$ws.Range("I2").Value = 3.69
$ws.Range("J2").Value = 2.405
$ws.Range("K2").Value = 5.875

# --- Row 3: Carbapenems ---
$ws.Range("D3").Value = 60
$ws.Range("E3").Value = 10423
$ws.Range("F3").Value = 22803
$ws.Range("G3").Value = 3661
$ws.Range("H3").Value = 2531
$ws.Range("J3").Value = 3.81
$ws.Range("K3").Value = 9.574999999999999

# --- Row 4: Cephems ---
$ws.Range("D4").Value = 44
$ws.Range("E4").Value = 11492
$ws.Range("F4").Value = 25037
$ws.Range("G4").Value = 3006
$ws.Range("H4").Value = 2873
$ws.Range("J4").Value = 2.42
$ws.Range("K4").Value = 5.84

# --- Row 5 ---
$ws.Range("D5").Value = 50
$ws.Range("E5").Value = 8321
$ws.Range("F5").Value = 18939
$ws.Range("G5").Value = 1797
$ws.Range("H5").Value = 1826
$ws.Range("I5").Value = 3.78
$ws.Range("J5").Value = 2.77
$ws.Range("K5").Value = 7.44

# --- Row 6 ---
$ws.Range("D6").Value = 31
$ws.Range("E6").Value = 4495
$ws.Range("F6").Value = 8310
$ws.Range("G6").Value = 1049
$ws.Range("H6").Value = 660
$ws.Range("I6").Value = 5.21
$ws.Range("K6").Value = 10.065

# --- Row 9 ---
$ws.Range("D9").Value = 125
$ws.Range("E9").Value = 48755
$ws.Range("F9").Value = 111634
$ws.Range("G9").Value = 17042
$ws.Range("H9").Value = 22826
$ws.Range("I9").Value = 3.41
$ws.Range("J9").Value = 1.8575
$ws.Range("K9").Value = 6.03

# --- Row 10 ---
$ws.Range("D10").Value = 37
$ws.Range("E10").Value = 4878
$ws.Range("F10").Value = 9611
$ws.Range("G10").Value = 2232
$ws.Range("H10").Value = 993
$ws.Range("I10").Value = 12.6
$ws.Range("J10").Value = 4.455
$ws.Range("K10").Value = 43.15

# --- Row 40: Other ---
$ws.Range("D40").Value = 111
$ws.Range("E40").Value = 77678
$ws.Range("F40").Value = 213465
$ws.Range("G40").Value = 17497.5
$ws.Range("H40").Value = 59716.2
$ws.Range("I40").Value = 2.72
$ws.Range("J40").Value = 0.585
$ws.Range("K40").Value = 4.73

# --- Remove the duplicate "Primary/specific infection site" row (old row 47). ---
# This shifts old rows 48-55 up to become new rows 47-54.
$ws.Rows("47").Delete()

# --- New row 49 ("Urinary tract") gets corrected figures. ---
$ws.Range("D49").Value = 16
$ws.Range("E49").Value = 3174
$ws.Range("F49").Value = 6563
$ws.Range("G49").Value = 567
$ws.Range("H49").Value = 1150
$ws.Range("I49").Value = 1.215
$ws.Range("J49").Value = 0.3025
$ws.Range("K49").Value = 2.4275

# --- New row 54 ("Prior colonization or infection") gains an A-column label. ---
$ws.Range("A54").Value = "Prior colonization or infection"

Write-Host "Edit applied"
